$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Through 2022-03-13" to "Through 2022-03-14"
$ws.Name = "Through 2022-03-14"

# Update the header label for the "current" month column (B)
$ws.Range("B1").Value = "March 2022 (through March 14)"

# Austin (row 3): B3 5 -> 6, K3 2 -> 3
$ws.Range("B3").Value = 6
$ws.Range("K3").Value = 3

# North Lawndale (row 4): E4 4 -> 5
$ws.Range("E4").Value = 5

# Garfield Park (row 5): W5 1 -> 2
$ws.Range("W5").Value = 2

# Little Italy, UIC (row 16): new E16 = 1
$ws.Range("E16").Value = 1

# Lincoln Park (row 19): new B19 = 1
$ws.Range("B19").Value = 1

# Chatham (row 22): new B22 = 1
$ws.Range("B22").Value = 1

# Riverdale (row 80): new N80 = 1
$ws.Range("N80").Value = 1
